# Update the "Generated At" timestamp on the Summary sheet.
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "2026-01-17 20:26 UTC"

# Fill in Market Cap (col E) and 24h Volume (col F) for rows 2-11 on each
# of the three setup sheets, replacing the "N/A" placeholders.

$reversalData = @(
    @("$740.18M", "$1.95M"),
    @("$394.28M", "$2.45M"),
    @("$597.12M", "$1.13M"),
    @("$365.91M", "$1.48M"),
    @("$1.15B", "$4.57M"),
    @("$753.93M", "$3.47M"),
    @("$1.26B", "$2.71M"),
    @("$505.02M", "$1.15M"),
    @("$384.17M", "$1.07M"),
    @("$2.69B", "$5.98M")
)

$breakoutData = @(
    @("$334.63M", "$14.23M"),
    @("$131.40M", "$4.31M"),
    @("$394.28M", "$2.45M"),
    @("$981.15M", "$2.22M"),
    @("$145.43M", "$3.40M"),
    @("$955.28M", "$38.78M"),
    @("$647.58M", "$18.71M"),
    @("$2.18B", "$9.79M"),
    @("$200.51M", "$8.78M"),
    @("$1.03B", "$7.74M")
)

$pullbackData = @(
    @("$334.63M", "$14.23M"),
    @("$131.40M", "$4.31M"),
    @("$394.28M", "$2.45M"),
    @("$647.58M", "$18.71M"),
    @("$981.15M", "$2.22M"),
    @("$145.43M", "$3.40M"),
    @("$121.99M", "$1.78M"),
    @("$2.92B", "$3.65M"),
    @("$664.65M", "$1.09M"),
    @("$365.91M", "$1.48M")
)

$sheetsToUpdate = @(
    @("Reversal Setups", $reversalData),
    @("Breakout Setups", $breakoutData),
    @("Pullback Setups", $pullbackData)
)

foreach ($entry in $sheetsToUpdate) {
    $sheetName = $entry[0]
    $data = $entry[1]
    $ws = $wb.Worksheets.Item($sheetName)

    for ($i = 0; $i -lt $data.Length; $i++) {
        $row = $i + 2
        $marketCap = $data[$i][0]
        $volume = $data[$i][1]
        $ws.Cells.Item($row, 5).Value = $marketCap
        $ws.Cells.Item($row, 6).Value = $volume
    }
}
